$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 264, shifting existing rows 264:273 down to 265:274
$ws.Rows.Item(264).Insert()

# Populate the new row 264 with its data
$ws.Cells.Item(264, 1).Value = 4
$ws.Cells.Item(264, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(264, 3).Value = "Los Lagos"
$ws.Cells.Item(264, 4).Value = 44509
$ws.Cells.Item(264, 5).Value = 10
$ws.Cells.Item(264, 6).Value = "Fruta"
$ws.Cells.Item(264, 7).Value = 100108
$ws.Cells.Item(264, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(264, 9).Value = 100108006
$ws.Cells.Item(264, 10).Value = "Plátano"
$ws.Cells.Item(264, 11).Value = "Sin especificar"
$ws.Cells.Item(264, 12).Value = "Primera Pintón"
$ws.Cells.Item(264, 13).Value = 1400
$ws.Cells.Item(264, 14).Value = 2300
$ws.Cells.Item(264, 15).Value = 22500
$ws.Cells.Item(264, 16).Value = 12400
$ws.Cells.Item(264, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(264, 18).Value = "Ecuador"
$ws.Cells.Item(264, 19).Value = 620
$ws.Cells.Item(264, 20).Value = 20
